$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.158.82'
$ws.Range('E2').Value = '  -0.10%  '

$ws.Range('D3').Value = '1.899.73'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.31%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.57'
$ws.Range('E5').Value = '  -0.58%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.26%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4632'
$ws.Range('E7').Value = '  -0.17%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3909'
$ws.Range('E8').Value = '  -1.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07879'
$ws.Range('E9').Value = '  -1.00%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9901'
$ws.Range('E10').Value = '  -0.92%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.81'
$ws.Range('E11').Value = '  -1.96%  '

$ws.Range('D12').Value = '1.877.31'
$ws.Range('E12').Value = '  -3.46%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.059'
$ws.Range('E13').Value = '  -0.76%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.742'
$ws.Range('E14').Value = '  -0.31%  '

$ws.Range('E15').Value = '  +0.75%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.12'
$ws.Range('E16').Value = '  -0.62%  '

$ws.Range('E17').Value = '  -0.23%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009984'
$ws.Range('E18').Value = '  -0.85%  '

$ws.Range('E19').Value = '  -0.50%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.23%  '

$ws.Range('D21').Value = '29.170.39'
$ws.Range('E21').Value = '  -0.21%  '

$ws.Range('E22').Value = '  -1.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.07'
$ws.Range('E23').Value = '  -0.08%  '

$ws.Range('D24').Value = '2.122.26'
$ws.Range('E24').Value = '  -1.78%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.109'
$ws.Range('E25').Value = '  +2.75%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.91'
$ws.Range('E26').Value = '  -0.53%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.44'
$ws.Range('E27').Value = '  -0.37%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.952'
$ws.Range('E28').Value = '  +0.55%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '118.31'
$ws.Range('E29').Value = '  -0.66%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.882'
$ws.Range('E30').Value = '  -5.65%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09339'
$ws.Range('E31').Value = '  -0.47%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9012'
$ws.Range('E32').Value = '  -2.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.249'
$ws.Range('E33').Value = '  -1.82%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.326'
$ws.Range('E34').Value = '  -1.56%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.174'
$ws.Range('E35').Value = '  -2.75%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.180'
$ws.Range('E36').Value = '  +0.24%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05765'
$ws.Range('E37').Value = '  -1.04%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02083'
$ws.Range('E38').Value = '  -0.94%  '

$ws.Range('E39').Value = '  -0.18%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.720'
$ws.Range('E40').Value = '  -3.33%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5698'
$ws.Range('E41').Value = '  -0.89%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1790'
$ws.Range('E42').Value = '  -0.74%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.711'
$ws.Range('E43').Value = '  -2.49%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.90'
$ws.Range('E44').Value = '  -0.98%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5347'
$ws.Range('E45').Value = '  -1.39%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.189'
$ws.Range('E46').Value = '  -0.91%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07010'
$ws.Range('E47').Value = '  -1.30%  '

$ws.Range('E48').Value = '  -1.29%  '

$ws.Range('E49').Value = '  -0.62%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '113.05'
$ws.Range('E50').Value = '  +0.80%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.048'
$ws.Range('E51').Value = '  -1.23%  '
